# Apply updated Betfair Back/Lay odds values for 2025-12-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 5.5
$ws.Range("H5").Value = 1.68
$ws.Range("I5").Value = 1.75
$ws.Range("K5").Value = 4.7
$ws.Range("Q5").Value = 1.52
# Row 6
$ws.Range("F6").Value = 2.26
$ws.Range("I6").Value = 4
$ws.Range("P6").Value = 1.63
$ws.Range("Y6").Value = 11
$ws.Range("AH6").Value = 23
# Row 9
$ws.Range("H9").Value = 5.2
$ws.Range("J9").Value = 4
$ws.Range("P9").Value = 1.98
# Row 10
$ws.Range("P10").Value = 1.24
# Row 11
$ws.Range("F11").Value = 2.32
$ws.Range("H11").Value = 2.4
$ws.Range("J11").Value = 3.45
$ws.Range("K11").Value = 4.8
# Row 12
$ws.Range("Q12").Value = 1.62
# Row 14
$ws.Range("F14").Value = 2.42
$ws.Range("G14").Value = 2.92
$ws.Range("J14").Value = 2.96
$ws.Range("K14").Value = 3.4
$ws.Range("P14").Value = 1.55
$ws.Range("Q14").Value = 2.46
# Row 15
$ws.Range("J15").Value = 3.35
# Row 16
$ws.Range("F16").Value = 3.3
$ws.Range("H16").Value = 2.14
$ws.Range("P16").Value = 2.34
$ws.Range("Q16").Value = 1.62
# Row 17
$ws.Range("G17").Value = 3.7
$ws.Range("H17").Value = 2.14
# Row 18
$ws.Range("G18").Value = 3.85
$ws.Range("I18").Value = 2.36
$ws.Range("J18").Value = 3.5
# Row 19
$ws.Range("F19").Value = 2.1
$ws.Range("H19").Value = 3.25
$ws.Range("I19").Value = 3.55
$ws.Range("J19").Value = 4
$ws.Range("K19").Value = 4.5
# Row 22
$ws.Range("F22").Value = 5.3
$ws.Range("I22").Value = 1.73
$ws.Range("S22").Value = 3.05
$ws.Range("T22").Value = 1.82
# Row 23
$ws.Range("F23").Value = 3.05
$ws.Range("G23").Value = 3.4
$ws.Range("H23").Value = 2.38
$ws.Range("I23").Value = 2.6
$ws.Range("K23").Value = 3.65
$ws.Range("P23").Value = 1.83
$ws.Range("Q23").Value = 2.02
# Row 28
$ws.Range("F28").Value = 2.2
$ws.Range("G28").Value = 2.24
$ws.Range("H28").Value = 3.55
$ws.Range("N28").Value = 3.75
$ws.Range("P28").Value = 1.96
$ws.Range("AB28").Value = 11
$ws.Range("AG28").Value = 11
$ws.Range("AI28").Value = 55
# Row 29
$ws.Range("F29").Value = 3.35
$ws.Range("H29").Value = 2.26
$ws.Range("Q29").Value = 1.77
# Row 30
$ws.Range("F30").Value = 1.7
$ws.Range("G30").Value = 1.78
$ws.Range("H30").Value = 4.8
$ws.Range("J30").Value = 4.4
$ws.Range("P30").Value = 2.46
# Row 32
$ws.Range("F32").Value = 2.32
$ws.Range("G32").Value = 2.68
$ws.Range("H32").Value = 3.1
$ws.Range("I32").Value = 4.1
$ws.Range("J32").Value = 2.8
$ws.Range("K32").Value = 3.25
# Row 37
$ws.Range("F37").Value = 2.48
$ws.Range("G37").Value = 2.5
$ws.Range("K37").Value = 3.25
$ws.Range("P37").Value = 1.73
$ws.Range("Q37").Value = 2.3
$ws.Range("R37").Value = 1.28
# Row 38
$ws.Range("H38").Value = 3.4
$ws.Range("K38").Value = 3.75
# Row 41
$ws.Range("F41").Value = 2.04
$ws.Range("G41").Value = 2.08
$ws.Range("I41").Value = 4.3
$ws.Range("K41").Value = 3.65
$ws.Range("O41").Value = 1.41
$ws.Range("X41").Value = 13.5
$ws.Range("AL41").Value = 50
# Row 42
$ws.Range("F42").Value = 2.08
$ws.Range("G42").Value = 2.32
$ws.Range("H42").Value = 3.55
# Row 43
$ws.Range("F43").Value = 1.35
$ws.Range("G43").Value = 1.38
$ws.Range("H43").Value = 9.199999999999999
$ws.Range("I43").Value = 11
$ws.Range("P43").Value = 2.84
$ws.Range("Q43").Value = 1.41
# Row 44
$ws.Range("F44").Value = 2.86
$ws.Range("K44").Value = 3.4
# Row 45
$ws.Range("G45").Value = 3.7
$ws.Range("H45").Value = 2.26
$ws.Range("J45").Value = 3.3
$ws.Range("K45").Value = 3.6
# Row 46
$ws.Range("I46").Value = 1.91
# Row 48
$ws.Range("F48").Value = 2.08
$ws.Range("H48").Value = 2.12
$ws.Range("J48").Value = 2.32
# Row 49
$ws.Range("F49").Value = 3.3
$ws.Range("G49").Value = 3.7
$ws.Range("H49").Value = 2.3
$ws.Range("I49").Value = 2.46
$ws.Range("J49").Value = 3.35
# Row 53
$ws.Range("F53").Value = 4.2
$ws.Range("H53").Value = 1.86
# Row 54
$ws.Range("F54").Value = 2.86
# Row 56
$ws.Range("F56").Value = 2.3
$ws.Range("G56").Value = 2.44
$ws.Range("H56").Value = 3.15
$ws.Range("J56").Value = 3.2
$ws.Range("K56").Value = 3.7
$ws.Range("P56").Value = 1.51
$ws.Range("Q56").Value = 2.12
# Row 57
$ws.Range("F57").Value = 1.95
$ws.Range("K57").Value = 4.5
$ws.Range("P57").Value = 2.32
$ws.Range("Q57").Value = 1.61
# Row 58
$ws.Range("P58").Value = 2.06
$ws.Range("Q58").Value = 1.74
# Row 59
$ws.Range("F59").Value = 2.56
$ws.Range("H59").Value = 2.62
$ws.Range("I59").Value = 3.5
$ws.Range("J59").Value = 3.05
$ws.Range("Q59").Value = 1.95
# Row 60
$ws.Range("F60").Value = 2.16
$ws.Range("H60").Value = 3.6
$ws.Range("AD60").Value = 14.5
# Row 61
$ws.Range("F61").Value = 20
$ws.Range("G61").Value = 21
$ws.Range("H61").Value = 1.15
$ws.Range("I61").Value = 1.16
$ws.Range("R61").Value = 2.4
$ws.Range("S61").Value = 1.67
$ws.Range("T61").Value = 1.9
$ws.Range("U61").Value = 2.06
$ws.Range("AG61").Value = 1000
$ws.Range("AO61").Value = 2.4
# Row 67
$ws.Range("F67").Value = 2.32
$ws.Range("G67").Value = 2.46
$ws.Range("H67").Value = 3.6
$ws.Range("I67").Value = 3.7
$ws.Range("J67").Value = 3.15
$ws.Range("K67").Value = 3.45
$ws.Range("P67").Value = 1.69
$ws.Range("Q67").Value = 2.2
# Row 68
$ws.Range("T68").Value = 1.93
# Row 70
$ws.Range("P70").Value = 1.73
# Row 71
$ws.Range("N71").Value = 2.96
$ws.Range("S71").Value = 4.9
$ws.Range("AI71").Value = 70
$ws.Range("AJ71").Value = 40
# Row 72
$ws.Range("F72").Value = 2.1
$ws.Range("G72").Value = 2.28
$ws.Range("J72").Value = 2.92
# Row 73
$ws.Range("Q73").Value = 2.12
# Row 74
$ws.Range("F74").Value = 2.26
$ws.Range("H74").Value = 3.7
$ws.Range("I74").Value = 3.95
# Row 76
$ws.Range("F76").Value = 2
$ws.Range("I76").Value = 4.5
# Row 77
$ws.Range("F77").Value = 2.02
$ws.Range("I77").Value = 1.99
$ws.Range("P77").Value = 2.64
# Row 78
$ws.Range("H78").Value = 1.41
$ws.Range("I78").Value = 3.6
$ws.Range("K78").Value = 1000
$ws.Range("P78").Value = 1.9
$ws.Range("Q78").Value = 1.75
# Row 79
$ws.Range("G79").Value = 1.25
$ws.Range("H79").Value = 8
$ws.Range("I79").Value = 25
$ws.Range("J79").Value = 6.2
# Row 82
$ws.Range("F82").Value = 4.4
$ws.Range("G82").Value = 4.5
$ws.Range("H82").Value = 2.02
$ws.Range("I82").Value = 2.04
$ws.Range("N82").Value = 3.5
$ws.Range("AA82").Value = 28
$ws.Range("AC82").Value = 7.8
$ws.Range("AE82").Value = 24
$ws.Range("AL82").Value = 65
# Row 83
$ws.Range("O83").Value = 1.51
$ws.Range("P83").Value = 1.63
$ws.Range("Q83").Value = 2.5
$ws.Range("Y83").Value = 15.5
$ws.Range("Z83").Value = 50
$ws.Range("AA83").Value = 270
$ws.Range("AD83").Value = 26
$ws.Range("AI83").Value = 160
$ws.Range("AK83").Value = 24
$ws.Range("AN83").Value = 19.5
# Row 84
$ws.Range("F84").Value = 3.35
$ws.Range("I84").Value = 2.42
$ws.Range("P84").Value = 1.86
$ws.Range("Q84").Value = 2.02
# Row 85
$ws.Range("F85").Value = 1.81
$ws.Range("G85").Value = 1.95
$ws.Range("H85").Value = 5.1
$ws.Range("I85").Value = 5.8
$ws.Range("P85").Value = 1.61
# Row 86
$ws.Range("F86").Value = 1.78
$ws.Range("H86").Value = 4.1
# Row 88
$ws.Range("F88").Value = 1.84
$ws.Range("H88").Value = 4.3
$ws.Range("K88").Value = 4.1
